$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (S) by copying formatting from the existing
# "2021" column (R), then overwrite with the new figures.
$ws.Range("R3:R8").Copy($ws.Range("S3:S8"))

# New 2022 values
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 6.9031689452913012
$ws.Range("S5").Value = 44.306188104841333
$ws.Range("S6").Value = 318
$ws.Range("S7").Value = 2041
$ws.Range("S8").Value = 4606580

# Updated 2021 figures (R8 changed; R4/R5 become static computed values
# instead of formulas, following the same recalculation as before)
$ws.Range("R8").Value = 4513063
$ws.Range("R4").Value = 6.9132648934880807
$ws.Range("R5").Value = 42.321589572314856

# Update the active selection to match the author's final cursor position
$ws.Range("R13").Select() | Out-Null
